$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Plain Find/Replace (Range.Text / Find.Execute) causes this runtime to
# rebuild a touched paragraph's run list, silently dropping "structural"
# empty runs (<w:r/>) that sit next to the text run being edited whenever
# the neighbouring run has no distinguishing run-properties. Several
# paragraphs in this document legitimately carry such an empty <w:r/>
# (preserved, per the diff, as unrelated context) right before the text
# run we need to edit, so a naive Find/Replace would silently strip them.
#
# Instead, replace text by targeting the precise character Range of the
# old string and calling Range.InsertXML with a single-run OOXML package
# fragment (same "pkg:package" shape Range.WordOpenXML returns). InsertXML
# only overwrites the characters inside that Range, leaving sibling runs
# (and the owning paragraph's pPr) completely untouched.
# ---------------------------------------------------------------------------

function Replace-ExactText($old, $new, $rpr) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw "Text not found: $old"
    }
    $rng = $d.Range($idx, $idx + $old.Length)

    $esc = $new -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    $runInner = ""
    if ($rpr) {
        $runInner = "<w:rPr>" + $rpr + "</w:rPr>"
    }
    $runInner = $runInner + "<w:t xml:space=`"preserve`">" + $esc + "</w:t>"

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($pkg)
}

function Replace-AllExactText($old, $new, $rpr) {
    while ($true) {
        $full = $d.Content.Text
        $idx = $full.IndexOf($old)
        if ($idx -lt 0) { break }
        Replace-ExactText $old $new $rpr
    }
}

# Replace the long, more specific sentence first so it does not get
# partially mangled by the shorter bullet-point replacement below (the
# bullet text is a substring of this sentence).
Replace-ExactText "Play Lightning Joker for free and enjoy exciting bonus features such as Respins and the Mystery Win feature. High payouts of up to 10,000 times the bet value." "Read our review of Lightning Joker, a modern online slot game with exciting features and high payouts. Play for free now!" "<w:i/>"

Replace-ExactText "Exciting bonus features, including Respins and a Mystery Win feature" "Exciting bonus features with respins and random multipliers" $null

Replace-ExactText "High payouts of up to 10,000 times the bet value" "High maximum payout of 10,000 times the bet" $null

Replace-ExactText "Classic fruit theme with a modern, elegant design" "Modern and elegant design with classic fruit symbols" $null

Replace-ExactText "Multiplier that reaches up to 10x" "Standard RTP of 96.3% for the best online gambling titles" $null

Replace-ExactText "Only 3 reels and 5 paylines may not appeal to players looking for more complex gameplay" "Limited number of paylines with only 5 available" $null

Replace-ExactText "The bonus features can be difficult to activate" "Lack of immersive storyline or theme" $null

# This title/heading text occurs twice: once as the plain Heading1 run,
# once as a bold run near the bottom. Replace each occurrence with the
# correct formatting explicitly, since InsertXML overwrites whatever run
# properties we give it for the targeted range.
Replace-ExactText "Play Lightning Joker for Free - Exciting Bonus Features" "Play Lightning Joker Free: Exciting Features and High Payouts" $null
Replace-ExactText "Play Lightning Joker for Free - Exciting Bonus Features" "Play Lightning Joker Free: Exciting Features and High Payouts" "<w:b/>"
